$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# New forecast values (Amazon Mean / P70 / P80 / P90) per week row, 2..17
# Columns: D = Amazon Mean Forecast, E = Amazon P70 Forecast,
#          F = Amazon P80 Forecast, G = Amazon P90 Forecast
$values = @{
    2  = @(64, 76, 87, 103)
    3  = @(60, 71, 82, 98)
    4  = @(61, 72, 82, 98)
    5  = @(60, 72, 83, 99)
    6  = @(61, 73, 86, 105)
    7  = @(60, 72, 83, 102)
    8  = @(61, 74, 87, 107)
    9  = @(64, 77, 91, 113)
    10 = @(60, 72, 85, 104)
    11 = @(61, 74, 87, 109)
    12 = @(63, 77, 91, 114)
    13 = @(65, 79, 95, 120)
    14 = @(64, 77, 91, 114)
    15 = @(60, 73, 88, 113)
    16 = @(59, 72, 87, 110)
    17 = @(58, 70, 85, 108)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("D$row").Value = $rowVals[0]
    $ws.Range("E$row").Value = $rowVals[1]
    $ws.Range("F$row").Value = $rowVals[2]
    $ws.Range("G$row").Value = $rowVals[3]
}
